$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "Wine to discover (todo)" right after the
#    "European design. & varieties" sheet (this pushes the following sheets
#    -- Wine Storage & Service, Tastes & aromas of varieties, Wine & food
#    pairing, Varieties of designations -- one position further down).
# ---------------------------------------------------------------------------
$europeanDesign = $wb.Worksheets.Item("European design. & varieties")
$newSheet = $wb.Worksheets.Add($null, $europeanDesign)
$newSheet.Name = "Wine to discover (todo)"

$newSheet.Range("A1").Value = "French"
$newSheet.Range("B1").Value = 4981
$newSheet.Range("C1").Value = "Kindle"
$newSheet.Range("D1").Value = "Region / Appellation"
$newSheet.Range("D2").Select()

# ---------------------------------------------------------------------------
# 2. Append the new Italian denomination rows (41-55) to the
#    "European design. & varieties" sheet.
# ---------------------------------------------------------------------------
$rows = @(
    @("Barolo (Red), Piedmont, Italy", "Nebbiolo"),
    @("Barbaresco (Red), Piedmont, Italy", "Nebbiolo"),
    @("Gattinara (Red), Piedmont, Italy", "Nebbiolo, Bonarda"),
    @("Gavi, Piedmont (White)", "Cortese"),
    @("Roero arneis (White), Piedmont, Italy", "Arneis"),
    @("Chianti, Chianti classico (Red), Tuscany, Italy", "Sangiovese, Canaiolo and others"),
    @("Brunello di Montalcino (Red), Tuscany, Italy", "Sangiovese"),
    @("Vernaccia di San Gimignano (White), Tuscany, Italy", "Vernaccia"),
    @("Vino noble di Montepulciano (Red), Tuscany, Italy", "Sangiovese, Canaiolo and others"),
    @("Carmignano (Red), Tuscany, Italy", "Sangiovese, Cabernet Sauvignon, Canaiolo and others"),
    @("Super-toscans (Red), Tuscany, Italy", " Cabernet Sauvignon, Sangiovese, and others"),
    @("Soave (White), Veneto, Italy", "Garganega, Trebbiano et autres"),
    @("Amarone (Red), Veneto, Italy", "Corvina, Molinara, Rondinella"),
    @("Bianco di Custoza (White), Veneto, Italy", "Trebbiano, Garganega, Tocai"),
    @("Luagana (White), Veneto, Italy", "Trebbiano")
)

$startRow = 41
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $europeanDesign.Cells.Item($r, 1).Value = $rows[$i][0]
    $europeanDesign.Cells.Item($r, 2).Value = $rows[$i][1]
}

$europeanDesign.Range("B55").Select()
$europeanDesign.Activate()

Write-Output "Italian denominations added"
